# Fruta / hortaliza, semanal
# The data rows (2-19) get their per-record columns (D, L:T) permuted
# across rows while the lookup columns (A,B,C,E:K) stay identical.
# Equivalent to reordering whole data rows according to the mapping below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: new row number -> old row number whose data it should receive
$map = @{
    2  = 8
    3  = 3
    4  = 18
    5  = 15
    6  = 6
    7  = 2
    8  = 7
    9  = 4
    10 = 11
    11 = 12
    12 = 10
    13 = 14
    14 = 13
    15 = 17
    16 = 19
    17 = 16
    18 = 5
    19 = 9
}

# Columns that actually carry the per-record values that get shuffled
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)   # D, L, M, N, O, P, Q, R, S, T

# Snapshot the original values for every relevant cell before overwriting anything
$snapshot = @{}
for ($r = 2; $r -le 19; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r`_$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Write back the values according to the permutation map
foreach ($newRow in $map.Keys) {
    $oldRow = $map[$newRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($newRow, $c).Value2 = $snapshot["$oldRow`_$c"]
    }
}
